$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kazbegi")
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1875
$ws.Range("K5").Value = 851
$ws.Range("K6").Value = 1024
